$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bioSampleNumber column (C) for rows 2-19: add 136 to existing value,
# and tighten row height to 13.8 for each of those rows.
for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 136
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Move the active selection to H8, matching the saved selection state.
$ws.Range("H8").Select()
